$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('B2') 'Bitcoin'
Set-TextValue $ws.Range('C2') 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextValue $ws.Range('D2') '29.354.48'
Set-TextValue $ws.Range('E2') '  +0.82%  '

Set-TextValue $ws.Range('B3') 'Ethereum'
Set-TextValue $ws.Range('C3') 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextValue $ws.Range('D3') '1.944.36'
Set-TextValue $ws.Range('E3') '  +2.37%  '

Set-TextValue $ws.Range('B4') 'TetherUSD'
Set-TextValue $ws.Range('C4') 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue $ws.Range('D4') '1.009'
Set-TextValue $ws.Range('E4') '  +0.72%  '

Set-TextValue $ws.Range('B5') 'BNB'
Set-TextValue $ws.Range('C5') 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range('D5') '325.77'
Set-TextValue $ws.Range('E5') '  +0.13%  '

Set-TextValue $ws.Range('B6') 'USDC'
Set-TextValue $ws.Range('C6') 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range('D6') '1.005'
Set-TextValue $ws.Range('E6') '  +0.34%  '

Set-TextValue $ws.Range('B7') 'XRP'
Set-TextValue $ws.Range('C7') 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range('D7') '0.4627'
Set-TextValue $ws.Range('E7') '  +0.33%  '

Set-TextValue $ws.Range('B8') 'Cardano'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D8') '0.3868'
Set-TextValue $ws.Range('E8') '  -0.76%  '

Set-TextValue $ws.Range('B9') 'OKB'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D9') '46.07'
Set-TextValue $ws.Range('E9') '  +0.46%  '

Set-TextValue $ws.Range('B10') 'Dogecoin'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D10') '0.07839'
Set-TextValue $ws.Range('E10') '  -0.59%  '

Set-TextValue $ws.Range('B11') 'Polygon'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D11') '0.9800'
Set-TextValue $ws.Range('E11') '  -1.07%  '

Set-TextValue $ws.Range('B12') 'Solana'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range('D12') '22.62'
Set-TextValue $ws.Range('E12') '  +3.52%  '

Set-TextValue $ws.Range('B13') 'WrappedEther'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.961.83'
Set-TextValue $ws.Range('E13') '  +3.74%  '

Set-TextValue $ws.Range('B14') 'Chainlink'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D14') '7.088'
Set-TextValue $ws.Range('E14') '  +0.41%  '

Set-TextValue $ws.Range('B15') 'Polkadot'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D15') '5.750'
Set-TextValue $ws.Range('E15') '  -0.31%  '

Set-TextValue $ws.Range('B16') 'TRON'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D16') '0.07067'
Set-TextValue $ws.Range('E16') '  +1.12%  '

Set-TextValue $ws.Range('B17') 'Litecoin'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D17') '86.91'
Set-TextValue $ws.Range('E17') '  -1.18%  '

Set-TextValue $ws.Range('B18') 'BinanceUSD'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D18') '1.010'
Set-TextValue $ws.Range('E18') '  +0.81%  '

Set-TextValue $ws.Range('B19') 'ShibaInu'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D19') '0.000009837'
Set-TextValue $ws.Range('E19') '  -1.38%  '

Set-TextValue $ws.Range('B20') 'Avalanche'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D20') '16.98'
Set-TextValue $ws.Range('E20') '  -0.52%  '

Set-TextValue $ws.Range('B21') 'Dai'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D21') '1.002'
Set-TextValue $ws.Range('E21') '  +0.12%  '

Set-TextValue $ws.Range('B22') 'WrappedBTC'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D22') '29.406.49'
Set-TextValue $ws.Range('E22') '  +0.93%  '

Set-TextValue $ws.Range('B23') 'Uniswap'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D23') '5.500'
Set-TextValue $ws.Range('E23') '  +3.39%  '

Set-TextValue $ws.Range('B24') 'Cosmos'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D24') '11.07'
Set-TextValue $ws.Range('E24') '  -0.51%  '

Set-TextValue $ws.Range('B25') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D25') '2.199.32'
Set-TextValue $ws.Range('E25') '  +3.89%  '

Set-TextValue $ws.Range('B26') 'Toncoin'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D26') '2.107'
Set-TextValue $ws.Range('E26') '  -0.21%  '

Set-TextValue $ws.Range('B27') 'Monero'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D27') '157.61'
Set-TextValue $ws.Range('E27') '  +1.06%  '

Set-TextValue $ws.Range('B28') 'EthereumClassic'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D28') '19.39'
Set-TextValue $ws.Range('E28') '  -0.33%  '

Set-TextValue $ws.Range('B29') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D29') '5.710'
Set-TextValue $ws.Range('E29') '  -3.37%  '

Set-TextValue $ws.Range('B30') 'BitcoinCash'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D30') '118.82'
Set-TextValue $ws.Range('E30') '  +0.40%  '

Set-TextValue $ws.Range('B31') 'LidoDAOToken'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D31') '1.850'
Set-TextValue $ws.Range('E31') '  -1.63%  '

Set-TextValue $ws.Range('B32') 'Stellar'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D32') '0.09377'
Set-TextValue $ws.Range('E32') '  +0.46%  '

Set-TextValue $ws.Range('B33') 'ImmutableX'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D33') '0.8590'
Set-TextValue $ws.Range('E33') '  -4.44%  '

Set-TextValue $ws.Range('B34') 'Filecoin'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D34') '5.172'
Set-TextValue $ws.Range('E34') '  -1.55%  '

Set-TextValue $ws.Range('B35') 'ARBITRUM'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D35') '1.306'
Set-TextValue $ws.Range('E35') '  -1.43%  '

Set-TextValue $ws.Range('B36') 'HuobiToken'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D36') '3.136'
Set-TextValue $ws.Range('E36') '  -0.66%  '

Set-TextValue $ws.Range('B37') 'Hedera'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D37') '0.05751'
Set-TextValue $ws.Range('E37') '  -0.66%  '

Set-TextValue $ws.Range('B38') 'TrustWalletToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D38') '1.152'
Set-TextValue $ws.Range('E38') '  -1.93%  '

Set-TextValue $ws.Range('B39') 'VeChain'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D39') '0.02080'
Set-TextValue $ws.Range('E39') '  -0.26%  '

Set-TextValue $ws.Range('B40') 'FraxShare'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D40') '7.632'
Set-TextValue $ws.Range('E40') '  -1.18%  '

Set-TextValue $ws.Range('B41') 'TheSandbox'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D41') '0.5652'
Set-TextValue $ws.Range('E41') '  -0.69%  '

Set-TextValue $ws.Range('B42') 'Algorand'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D42') '0.1779'
Set-TextValue $ws.Range('E42') '  -0.68%  '

Set-TextValue $ws.Range('B43') 'Aptos'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D43') '9.390'
Set-TextValue $ws.Range('E43') '  -3.35%  '

Set-TextValue $ws.Range('B44') 'MXToken'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D44') '2.730'
Set-TextValue $ws.Range('E44') '  +6.86%  '

Set-TextValue $ws.Range('B45') 'PEPE'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D45') '0.000002776'
Set-TextValue $ws.Range('E45') '  +43.04%  '

Set-TextValue $ws.Range('B46') 'EnergySwap'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '11.62'
Set-TextValue $ws.Range('E46') '  -2.73%  '

Set-TextValue $ws.Range('B47') 'Decentraland'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D47') '0.5277'
Set-TextValue $ws.Range('E47') '  -1.44%  '

Set-TextValue $ws.Range('B48') 'RenderToken'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D48') '2.111'
Set-TextValue $ws.Range('E48') '  -5.49%  '

Set-TextValue $ws.Range('B49') 'Cronos'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D49') '0.06901'
Set-TextValue $ws.Range('E49') '  -1.58%  '

Set-TextValue $ws.Range('B50') 'NEARProtocol'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D50') '1.807'
Set-TextValue $ws.Range('E50') '  -2.26%  '

Set-TextValue $ws.Range('B51') 'Quant'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D51') '111.59'
Set-TextValue $ws.Range('E51') '  -1.30%  '
